$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("F4").Value = 1.55
$ws.Range("G4").Value = 1.61
$ws.Range("H4").Value = 6.4
$ws.Range("I4").Value = 7.8
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 5.1
$ws.Range("P4").Value = 1.99

# Row 10 update
$ws.Range("L10").Value = 1.47
